# Apply the stock-report corrections described in the commit diff.
# Each line sets a single cell on the (only) worksheet to its new value.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F93").Value = 13
$ws.Range("G93").Value = 440.05
$ws.Range("B124").Value = 416546.71
$ws.Range("B161").Value = 53925
$ws.Range("F161").Value = 1
$ws.Range("G161").Value = 66.44
$ws.Range("B162").Value = 64350
$ws.Range("E162").Value = 70.63
$ws.Range("F162").Value = 92
$ws.Range("G162").Value = 6112.48
$ws.Range("B163").Value = 57756
$ws.Range("E163").Value = 79.37
$ws.Range("F163").Value = -100
$ws.Range("G163").Value = -6644
$ws.Range("B279").Value = 64973
$ws.Range("E279").Value = 35.4
$ws.Range("F279").Value = 146
$ws.Range("G279").Value = 4861.8
$ws.Range("B280").Value = 48706
$ws.Range("E280").Value = 39.8
$ws.Range("F280").Value = -144
$ws.Range("G280").Value = -4795.2
$ws.Range("B316").Value = 57077
$ws.Range("D316").Value = 93.08
$ws.Range("E316").Value = 111.2
$ws.Range("F316").Value = 1
$ws.Range("G316").Value = 93.08
$ws.Range("B317").Value = 61610
$ws.Range("D317").Value = 102.71
$ws.Range("E317").Value = 122.71
$ws.Range("F317").Value = -58
$ws.Range("G317").Value = -5957.18
$ws.Range("B318").Value = 63565
$ws.Range("E318").Value = 109.19
$ws.Range("F318").Value = 60
$ws.Range("G318").Value = 6162.6
$ws.Range("F319").Value = 102
$ws.Range("G319").Value = 6482.1
$ws.Range("F324").Value = 15
$ws.Range("G324").Value = 1244.1
$ws.Range("F328").Value = 257
$ws.Range("G328").Value = 41351.3
$ws.Range("B346").Value = 63520
$ws.Range("E346").Value = 153.4
$ws.Range("F346").Value = 95
$ws.Range("G346").Value = 13706.6
$ws.Range("B347").Value = 55373
$ws.Range("E347").Value = 163.62
$ws.Range("F347").Value = -94
$ws.Range("G347").Value = -13562.32
$ws.Range("B350").Value = 63531
$ws.Range("E350").Value = 152.53
$ws.Range("F350").Value = 80
$ws.Range("G350").Value = 11478.4
$ws.Range("B352").Value = 57802
$ws.Range("E352").Value = 162.71
$ws.Range("F352").Value = -79
$ws.Range("G352").Value = -11334.92
$ws.Range("B372").Value = 63652
$ws.Range("E372").Value = 55.42
$ws.Range("F372").Value = 243
$ws.Range("G372").Value = 12667.59
$ws.Range("B373").Value = 57885
$ws.Range("E373").Value = 62.28
$ws.Range("F373").Value = 4
$ws.Range("G373").Value = 208.52
$ws.Range("B375").Value = 61605
$ws.Range("E375").Value = 133.78
$ws.Range("F375").Value = -13
$ws.Range("G375").Value = -1455.48
$ws.Range("B376").Value = 63563
$ws.Range("E376").Value = 119.04
$ws.Range("F376").Value = 15
$ws.Range("G376").Value = 1679.4
$ws.Range("B379").Value = 61608
$ws.Range("E379").Value = 154.12
$ws.Range("F379").Value = -56
$ws.Range("G379").Value = -7224.56
$ws.Range("B380").Value = 63564
$ws.Range("E380").Value = 137.16
$ws.Range("F380").Value = 57
$ws.Range("G380").Value = 7353.57
$ws.Range("B389").Value = 62865
$ws.Range("F389").Value = 129
$ws.Range("G389").Value = 10295.49
$ws.Range("B390").Value = 57817
$ws.Range("F390").Value = 3
$ws.Range("G390").Value = 239.43
$ws.Range("B419").Value = 63007
$ws.Range("F419").Value = 952
$ws.Range("G419").Value = 163106.16
$ws.Range("B420").Value = 57856
$ws.Range("F420").Value = 2
$ws.Range("G420").Value = 342.66
$ws.Range("B431").Value = 53082
$ws.Range("C431").Value = "HUL-VIM BAR MULTIPACK FW 4X200G"
$ws.Range("F431").Value = 1
$ws.Range("G431").Value = 59.47
$ws.Range("B432").Value = 63102
$ws.Range("C432").Value = "HUL-Vim Bar Multipack Fw 4X200G"
$ws.Range("F432").Value = 4
$ws.Range("G432").Value = 237.88
$ws.Range("B434").Value = 731586.7
$ws.Range("B457").Value = 63681
$ws.Range("E457").Value = 23.84
$ws.Range("F457").Value = 56
$ws.Range("G457").Value = 1255.52
$ws.Range("B458").Value = 31930
$ws.Range("E458").Value = 26.8
$ws.Range("F458").Value = -62
$ws.Range("G458").Value = -1390.04
$ws.Range("F502").Value = 161
$ws.Range("G502").Value = 33953.29
$ws.Range("B504").Value = 356501.69
$ws.Range("F525").Value = 89
$ws.Range("G525").Value = 13374.03
$ws.Range("B526").Value = 31522.94
$ws.Range("B583").Value = 65066
$ws.Range("E583").Value = 13.61
$ws.Range("F583").Value = 297
$ws.Range("G583").Value = 3804.57
$ws.Range("B584").Value = 53263
$ws.Range("E584").Value = 15.29
$ws.Range("F584").Value = -309
$ws.Range("G584").Value = -3958.29
$ws.Range("B586").Value = 45695
$ws.Range("E586").Value = 23.58
$ws.Range("F586").Value = -36
$ws.Range("G586").Value = -710.28
$ws.Range("B587").Value = 64915
$ws.Range("E587").Value = 20.98
$ws.Range("F587").Value = 33
$ws.Range("G587").Value = 651.09
$ws.Range("B599").Value = 64925
$ws.Range("E599").Value = 13.97
$ws.Range("F599").Value = 296
$ws.Range("G599").Value = 3892.4
$ws.Range("B600").Value = 45709
$ws.Range("E600").Value = 15.69
$ws.Range("F600").Value = -300
$ws.Range("G600").Value = -3945
$ws.Range("F604").Value = 325
$ws.Range("G604").Value = 4787.25
$ws.Range("B606").Value = 19754.6
$ws.Range("F629").Value = 362
$ws.Range("G629").Value = 8778.5
$ws.Range("B651").Value = 79490.39
$ws.Range("B709").Value = 64833
$ws.Range("E709").Value = 34.9
$ws.Range("F709").Value = 99
$ws.Range("G709").Value = 3250.17
$ws.Range("B710").Value = 60025
$ws.Range("E710").Value = 37.22
$ws.Range("F710").Value = -98
$ws.Range("G710").Value = -3217.34
$ws.Range("F843").Value = 88
$ws.Range("G843").Value = 11398.64
$ws.Range("F848").Value = 24
$ws.Range("G848").Value = 8160.72
$ws.Range("B849").Value = 102052.35
$ws.Range("B952").Value = 5635279.73
$ws.Range("B953").Value = 5635279.73
